$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cartesian Velocity - Linear")

# Row 3 - Standard Deviation
$ws.Range("B3").Value = 0.00383111640114036
$ws.Range("C3").Value = 0.006255949698204192
$ws.Range("D3").Value = 0.01945997593335323
$ws.Range("E3").Value = 0.01022964780240761

# Row 4 - Maximum
$ws.Range("B4").Value = 0.01663162550035524
$ws.Range("C4").Value = 0.02773520515886496
$ws.Range("D4").Value = 0.07739125956626192
$ws.Range("E4").Value = 0.06413217448647635

# Row 5 - Mean
$ws.Range("B5").Value = 0.006196401264470401
$ws.Range("C5").Value = 0.009455005521889643
$ws.Range("D5").Value = 0.0263313692048925
$ws.Range("E5").Value = 0.01599737928832565
